# Textbox response formatting fix
# Rename sheets and update stimulus filenames per updated task-order timestamps.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-16511687241706667"
$ws1.Range("B2").Value = "go_stims-1651168724138656.csv"
$ws1.Range("B3").Value = "GNG_stims-16511687241534302.csv"
$ws1.Range("B4").Value = "go_stims-16511687241561973.csv"
$ws1.Range("B5").Value = "GNG_stims-16511687241686034.csv"

# --- Sheet 2: NB ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-16511687273950837"
$ws2.Range("B2").Value = "OB-16511687259568806.csv"
$ws2.Range("B3").Value = "OB-16511687266421719.csv"
$ws2.Range("B4").Value = "ZB-match_0-16511687251179452.csv"
$ws2.Range("B5").Value = "TB-16511687266934159.csv"
$ws2.Range("B6").Value = "OB-16511687261811182.csv"
$ws2.Range("B7").Value = "TB-16511687273830278.csv"
$ws2.Range("B8").Value = "ZB-match_2-16511687259272697.csv"
$ws2.Range("B9").Value = "TB-16511687270310555.csv"
$ws2.Range("B10").Value = "ZB-match_4-16511687252759387.csv"

# --- Sheet 3: RS ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-1651168727396996"
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

# --- Sheet 4: TOL ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-16511687274586523"
$ws4.Range("B2").Value = "MM_stims-16511687274105837.csv"
$ws4.Range("B3").Value = "ZM_stims-1651168727398996.csv"
$ws4.Range("B4").Value = "MM_stims-16511687274417655.csv"
$ws4.Range("B5").Value = "ZM_stims-16511687274105837.csv"
$ws4.Range("B6").Value = "MM_stims-16511687274576838.csv"
$ws4.Range("B7").Value = "ZM_stims-16511687274427304.csv"

# --- Sheet 5: vSAT ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-16511687275351"
$ws5.Range("B2").Value = "SAT_stims-1651168727488369.csv"
$ws5.Range("B3").Value = "SAT_stims-16511687274620602.csv"
$ws5.Range("B4").Value = "vSAT_stims-1651168727503526.csv"
$ws5.Range("B5").Value = "vSAT_stims-165116872751962.csv"
